# "analyse sur les graphes" -- rename the data/analysis sheets and refresh
# the saved view state (selection / zoom) on the analysis sheet, matching
# the author's commit.

$wb = $excel.ActiveWorkbook

# --- Rename sheets: Sheet1 -> data_queens, Sheet2 -> analysis_queens -----
$wsData = $wb.Worksheets.Item("Sheet1")
$wsData.Name = "data_queens"

$wsAnalysis = $wb.Worksheets.Item("Sheet2")
$wsAnalysis.Name = "analysis_queens"

# --- Update the saved selection / zoom on the analysis sheet -------------
# analysis_queens is already the active sheet/tab in this workbook, so we
# can move the selection and set the zoom level without disturbing which
# tab is marked as selected.
[void]$wsAnalysis.Range("A20").Select()
$excel.ActiveWindow.Zoom = 100
